$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing sheets.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("MCEQs simplified").Name = "Macrochemical eq's simplified"
$wb.Worksheets.Item("lambdas").Name          = "Growth rates"
$wb.Worksheets.Item("MCEQs").Name            = "Macrochemical equations"

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "General info" sheet in front of everything else.
# ---------------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$infoSheet.Name = "General info"

# Plain intro line.
$infoSheet.Range("A1").Value = "This file contains source code for the following article:"

# --- A3: "Title: ..." with the word "Title" bold -----------------------
$infoSheet.Range("A3").Value = "Title: A Trade-off between Force and Flow may lead to Reduced Entropy Production Rate during Faster Microbial Growth"
$infoSheet.Range("A3").Characters(1, 5).Font.Bold = $true
$infoSheet.Range("A3").Characters(6, 111).Font.Bold = $false

# --- A4: "Authors: ..." with the word "Authors" bold --------------------
$infoSheet.Range("A4").Value = "Authors: Maarten J. Droste^(1,2), Maaike Remeijer^2, Robert Planqué^1, Frank J. Bruggeman^2"
$infoSheet.Range("A4").Characters(1, 7).Font.Bold = $true
$infoSheet.Range("A4").Characters(8, 84).Font.Bold = $false

# --- A5: "Affiliations: ..." with the word "Affiliations" bold ----------
$infoSheet.Range("A5").Value = "Affiliations: ^1Department of Mathematics, Amsterdam Center for Dynamics and Computation, Vrije Universiteit Amsterdam, 1081 HV Amsterdam, the Netherlands, ^2Systems Biology Lab, A-LIFE, AIMMS, Vrije Universiteit Amsterdam, 1081 HZ Amsterdam, the Netherlands"
$infoSheet.Range("A5").Characters(1, 12).Font.Bold = $true
$infoSheet.Range("A5").Characters(13, 246).Font.Bold = $false

# View settings for the new sheet.
$infoSheet.Activate()
$excel.ActiveWindow.Zoom = 120
[void]$infoSheet.Range("G14").Select()

# ---------------------------------------------------------------------------
# 3. Tweak the selection remembered on "Growth rates" (previously "lambdas").
# ---------------------------------------------------------------------------
$growthRates = $wb.Worksheets.Item("Growth rates")
$growthRates.Activate()
[void]$growthRates.Range("F22").Select()

# ---------------------------------------------------------------------------
# 4. Leave "Macrochemical eq's simplified" as the active / selected sheet
#    (its own remembered selection, AK9, is untouched by merely activating).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Macrochemical eq's simplified").Activate()

# ---------------------------------------------------------------------------
# 5. Cosmetic: rename the default cell style "Standaard" -> "Normal" (matches
#    the locale clean-up done when the file was re-saved from an EN build of
#    Excel).
# ---------------------------------------------------------------------------
foreach ($s in $wb.Styles) {
    if ($s.Name -eq "Standaard") {
        try { $s.NameLocal = "Normal" } catch {}
    }
}
